# Shradha AI Technologies Ltd_Semi_Final.xlsx - "Quarterly" sheet rework.
#
# The quarter-label column ("Jun 22 Q1", "Sep 22 Q2", ...) is split into
# three separate columns: Year (col A, kept), Month (new col B) and
# Quarter (new col C). All the former data columns (old B:U) shift two
# places to the right (new D:W).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Helper: write a value as TEXT (not auto-converted to a number) while
# leaving the cell with no explicit style/numberformat applied - matches
# how the rest of the sheet's data cells are stored (no "s" attribute).
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- 1. Insert two new, blank columns before the old column B -------------
# (old column B "Net sales/income from operations" etc. becomes column D)
$ws.Range("B1:C1").EntireColumn.Insert()

# --- 2. Header row ----------------------------------------------------------
$ws.Range("B1").Value = "Month"
$ws.Range("C1").Value = "Quarter"

# --- 3. Per-row Year / Month / Quarter values -------------------------------
$quarterInfo = @{
    2  = @("2022", "06", "Q1")
    3  = @("2022", "09", "Q2")
    4  = @("2022", "12", "Q3")
    5  = @("2023", "03", "Q4")
    6  = @("2023", "06", "Q1")
    7  = @("2023", "09", "Q2")
    8  = @("2023", "12", "Q3")
    9  = @("2024", "03", "Q4")
    10 = @("2024", "06", "Q1")
    11 = @("2024", "09", "Q2")
}

foreach ($row in 2..11) {
    $info = $quarterInfo[$row]
    $year = $info[0]
    $month = $info[1]
    $quarter = $info[2]

    Set-TextValue $ws.Range("A$row") $year
    Set-TextValue $ws.Range("B$row") $month
    $ws.Range("C$row").Value = $quarter
}
